# Update script by alec
# - Replace the "mahendra/activa" and "yogendra/scooty" rows with
#   "Ram/Cycle" and "Shyam/Tractor"
# - Remove the trailing "rahul/bullet" row entirely
# - Move the active selection to E11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Ram"
$ws.Range("B6").Value = "Cycle"
$ws.Range("A7").Value = "Shyam"
$ws.Range("B7").Value = "Tractor"

$ws.Rows("8:8").Delete()

$ws.Range("E11").Select()
